# Update "想去人数" (interested-count) figures for five events that are
# listed on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 2095
$ws1.Range("F9").Value = 10782
$ws1.Range("F15").Value = 9004
$ws1.Range("F18").Value = 5284
$ws1.Range("F20").Value = 3355

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2095
$ws4.Range("F12").Value = 10782
$ws4.Range("F18").Value = 9004
$ws4.Range("F21").Value = 5284
$ws4.Range("F23").Value = 3355
